$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column E (duplicate_image_filename) with "NA" for rows 2 through 21
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
